# edit.ps1
# Applies the Registru_pacienti.xlsx changes described in the commit diff:
#  - updates patient entry #1 (row 2) visit date + recommendation text
#  - reworks patient entry #3 (row 4) into the MBAPPEE KILLIAN record
#  - reworks patient entry #4 (row 5) into the PUFOSEL MOTAN record
#  - appends a brand-new patient entry #5 (row 6) for GIGI BECALI
#  - widens a few columns and extends dimension/autofilter to row 6

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width tweaks (col F/TELEFON, col I/JUDET, col V/RECOMANDARE) ---
# (ColumnWidth setter in this engine rounds through a +5/6 character-width
#  offset, so we subtract it up front to land exactly on the target widths)
$wOffset = 5/6
$ws.Columns(6).ColumnWidth = 22 - $wOffset
$ws.Columns(9).ColumnWidth = 21 - $wOffset
$ws.Columns(22).ColumnWidth = 200 - $wOffset

# --- Row 2: follow-up visit pushed back from 01-12-2023 to 20-12-2023 ---
$ws.Range("B2").Value = "20-12-2023"
$ws.Range("V2").Value = "A SE EVITA FRIGUL`nAERIUS LA NEVOIE`nA SE REPETA CONTROLUL DUPA 3 LUNI DE ZILE.`n"

# --- Row 4: KILLIAN MBAPPEE entry rewritten with new details ---
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "'13-12-2023"
$ws.Range("C4").Value = "KILLIAN"
$ws.Range("D4").Value = "MBAPPEE"
$ws.Range("E4").Value = "'1234567891011"
$ws.Range("F4").Value = "'078941256322"
$ws.Range("I4").Value = "IF-ILFOV"
$ws.Range("J4").Value = "YES"
$ws.Range("K4").Value = "Pensionar"
$ws.Range("L4").Value = "YES"
$ws.Range("M4").Value = "'25631214"
$ws.Range("N4").Value = "ASDSADSA`nASDSADSADSADSASDGFDBGFDDSGDFGDFGDF"
$ws.Range("O4").Value = "NO"
$ws.Range("P4").Value = "NON-APLICABIL"
$ws.Range("Q4").Value = "NON-APLICABIL"
$ws.Range("R4").Value = "NON-APLICABIL"
$ws.Range("S4").Value = "NON-APLICABIL"
$ws.Range("T4").Value = "YES"
$ws.Range("U4").Value = "EGOISM`nAPLICATIA VAD CA FUNCTIONEAZA ACUMA"
$ws.Range("V4").Value = "A SE MUTA LA REAL MADRID`nPOATE RAMANE SI LA PSG DACA VREA SA FACA RECORDURI, DAR NU CRED CA E BUN`nASASASA`nASASASQ`nMAESTRO KIMBPEMBPE EEEEE`nADADADADA`nADASDASDSADSA`nTEST TEST TEST TEST SCROLL`n"

# --- Row 5: PUFOSEL MOTAN entry rewritten with new details ---
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "'04-01-2024"
$ws.Range("C5").Value = "PUFOSEL"
$ws.Range("D5").Value = "MOTAN"
$ws.Range("E5").Value = "'1900106375492"
$ws.Range("F5").Value = "'0748313438"
$ws.Range("G5").Value = "PISICEASCA"
$ws.Range("H5").Value = "BUCOVINA"
$ws.Range("I5").Value = "SV-SUCEAVA"
$ws.Range("J5").Value = "YES"
$ws.Range("K5").Value = "Salariat"
$ws.Range("L5").Value = "YES"
$ws.Range("M5").Value = "XDQWDX"
$ws.Range("N5").Value = "PACIENT PISICOS, SE ADRESEAZA PENTRU RESPIRATOE ORALA NOCTURNA, ALINTATURI, MOTANELI.`n"
$ws.Range("O5").Value = "YES"
$ws.Range("P5").Value = "OBSTRUCTIV"
$ws.Range("Q5").Value = "NAZALA"
$ws.Range("R5").Value = "BUNA"
$ws.Range("S5").Value = "4-12"
$ws.Range("T5").Value = "YES"
$ws.Range("U5").Value = "HIPODRAGANEALA`n"
$ws.Range("V5").Value = "-CRESTE DRAGANEALA LA ZILNIC, 3 PUPURIX10/ZI, 20 IMBRATISARI/ZI, MINIM;`n-LA NEVOIE, SUPLIMENTERAZA DRAGANEALA.`n"

# --- Row 6: brand-new patient entry (GIGI BECALI) ---
# Copy the row-5 "ID" cell formatting (centered, bordered style) down to A6
# before filling in row 6 so the new row matches the look of the others.
$ws.Range("A5").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $False

# --- Row 6: new GIGI BECALI entry appended after row 5 ---
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "'20-12-2023"
$ws.Range("C6").Value = "GIGI"
$ws.Range("D6").Value = "BECALI"
$ws.Range("E6").Value = "'1235558971414"
$ws.Range("F6").Value = "'078945632"
$ws.Range("H6").Value = "PIPERA"
$ws.Range("I6").Value = "B-BUCURESTI"
$ws.Range("J6").Value = "YES"
$ws.Range("K6").Value = "Lege Speciala"
$ws.Range("L6").Value = "YES"
$ws.Range("M6").Value = "1213AAA"
$ws.Range("N6").Value = "PACIENT IN VARSTA DE 70 DE ANI SUFERA DE DEMENTA CHAMPIONS LEAGUE`n"
$ws.Range("O6").Value = "NO"
$ws.Range("P6").Value = "NON-APLICABIL"
$ws.Range("Q6").Value = "NON-APLICABIL"
$ws.Range("R6").Value = "NON-APLICABIL"
$ws.Range("S6").Value = "NON-APLICABIL"
$ws.Range("T6").Value = "YES"
$ws.Range("U6").Value = "OISM`nHAHALERISM`n"
$ws.Range("V6").Value = "A SE ASTEPTA PANA IN VARA ANULUI VIITOR CA POATE CASTIGA TITLUL. CONTROL IN IUNIE`n"

# --- Extend the table dimensions: dimension / AutoFilter / _FilterDatabase ---
$ws.AutoFilterMode = $False
$ws.Range("A1:V6").AutoFilter() | Out-Null
foreach ($n in $wb.Names) {
    if ($n.Name -eq "REGISTRU!_FilterDatabase") {
        $n.RefersTo = "='REGISTRU'!`$A`$1:`$V`$6"
    }
}
